$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '60.155.78'
$ws.Range('E2').Value = '  +0.12%  '
$ws.Range('D3').Value = '2.420.26'
$ws.Range('E3').Value = '  -0.10%  '
$ws.Range('E4').Value = '  -0.01%  '
$ws.Range('D5').Value = "'553.95"
$ws.Range('D5').Style = $ws.Range('B5').Style
$ws.Range('E5').Value = '  +0.06%  '
$ws.Range('D6').Value = "'137.37"
$ws.Range('D6').Style = $ws.Range('B6').Style
$ws.Range('E6').Value = '  -0.54%  '
$ws.Range('E7').Value = '  +0.03%  '
$ws.Range('D8').Value = "'0.589"
$ws.Range('D8').Style = $ws.Range('B8').Style
$ws.Range('E8').Value = '  +2.37%  '
$ws.Range('E9').Value = '  -1.50%  '
$ws.Range('D10').Value = "'5.72"
$ws.Range('D10').Style = $ws.Range('B10').Style
$ws.Range('E10').Value = '  -0.99%  '
$ws.Range('E11').Value = '  -0.23%  '
$ws.Range('E12').Value = '  -1.76%  '
$ws.Range('D13').Value = "'24.93"
$ws.Range('D13').Style = $ws.Range('B13').Style
$ws.Range('E13').Value = '  -0.09%  '
$ws.Range('D14').Value = '2.851.06'
$ws.Range('D15').Value = '60.064.16'
$ws.Range('E15').Value = '  +0.18%  '
$ws.Range('E16').Value = '  -0.55%  '
$ws.Range('D17').Value = '2.411.59'
$ws.Range('E17').Value = '  -0.14%  '
$ws.Range('E18').Value = '  -0.78%  '
$ws.Range('E19').Value = '  +2.47%  '
$ws.Range('D20').Value = "'327.41"
$ws.Range('D20').Style = $ws.Range('B20').Style
$ws.Range('E20').Value = '  -1.36%  '
$ws.Range('D21').Value = "'6.75"
$ws.Range('D21').Style = $ws.Range('B21').Style
$ws.Range('E21').Value = '  -0.28%  '
$ws.Range('E22').Value = '  +0.06%  '
$ws.Range('D23').Value = "'65.17"
$ws.Range('D23').Style = $ws.Range('B23').Style
$ws.Range('E23').Value = '  +0.11%  '
$ws.Range('E24').Value = '  +4.27%  '
$ws.Range('E25').Value = '  +1.44%  '
$ws.Range('E26').Value = '  +0.10%  '
$ws.Range('D27').Value = "'1.41"
$ws.Range('D27').Style = $ws.Range('B27').Style
$ws.Range('E27').Value = '  +5.13%  '
$ws.Range('E28').Value = '  -1.21%  '
$ws.Range('E29').Value = '  -0.44%  '
$ws.Range('D30').Value = "'170.29"
$ws.Range('D30').Style = $ws.Range('B30').Style
$ws.Range('E30').Value = '  +0.37%  '
$ws.Range('E31').Value = '  -2.67%  '
$ws.Range('D32').Value = "'0.405"
$ws.Range('D32').Style = $ws.Range('B32').Style
$ws.Range('E32').Value = '  -3.31%  '
$ws.Range('E33').Value = '  +2.46%  '
$ws.Range('D34').Value = "'18.55"
$ws.Range('D34').Style = $ws.Range('B34').Style
$ws.Range('E34').Value = '  -0.88%  '
$ws.Range('E35').Value = '  +0.02%  '
$ws.Range('E36').Value = '  +2.20%  '
$ws.Range('B37').Value = 'FirstDigitalUSD'
$ws.Range('C37').Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range('D37').Value = "'1.00"
$ws.Range('D37').Style = $ws.Range('B37').Style
$ws.Range('E37').Value = '  +0.03%  '
$ws.Range('B38').Value = 'NEARProtocol'
$ws.Range('C38').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D38').Value = "'4.22"
$ws.Range('D38').Style = $ws.Range('B38').Style
$ws.Range('E38').Value = '  +0.33%  '
$ws.Range('D39').Value = "'327.04"
$ws.Range('D39').Style = $ws.Range('B39').Style
$ws.Range('E39').Value = '  +2.55%  '
$ws.Range('E40').Value = '  -0.68%  '
$ws.Range('D41').Value = "'144.06"
$ws.Range('D41').Style = $ws.Range('B41').Style
$ws.Range('E41').Value = '  +3.20%  '
$ws.Range('E42').Value = '  -0.64%  '
$ws.Range('D43').Value = "'20.03"
$ws.Range('D43').Style = $ws.Range('B43').Style
$ws.Range('E43').Value = '  +2.24%  '
$ws.Range('E44').Value = '  +0.51%  '
$ws.Range('E45').Value = '  -0.61%  '
$ws.Range('E46').Value = '  +0.31%  '
$ws.Range('D47').Value = "'0.0223"
$ws.Range('D47').Style = $ws.Range('B47').Style
$ws.Range('E47').Value = '  -1.48%  '
$ws.Range('D48').Value = "'11.04"
$ws.Range('D48').Style = $ws.Range('B48').Style
$ws.Range('E48').Value = '  -0.03%  '
$ws.Range('E49').Value = '  -1.51%  '
$ws.Range('E50').Value = '  -0.45%  '
$ws.Range('E51').Value = '  -0.50%  '
